$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F for rows 2-25
$bf = New-Object "object[,]" 24,5
$bf[0,0] = 17.39238349376427
$bf[0,1] = 5.265214403894899
$bf[0,2] = 11.08490074222824
$bf[0,3] = 11.15595206546573
$bf[0,4] = 55.97864597782567
$bf[1,0] = 17.3439019773052
$bf[1,1] = 5.180421910722464
$bf[1,2] = 10.94327030893495
$bf[1,3] = 11.13865921087207
$bf[1,4] = 54.83860886638132
$bf[2,0] = 17.32078501600712
$bf[2,1] = 5.12605006425392
$bf[2,2] = 10.8543188840361
$bf[2,3] = 11.12945863325318
$bf[2,4] = 54.12850773851184
$bf[3,0] = 17.31304411219287
$bf[3,1] = 5.103313453015889
$bf[3,2] = 10.81758476478702
$bf[3,3] = 11.12606712026641
$bf[3,4] = 53.83687583918211
$bf[4,0] = 17.31186033411972
$bf[4,1] = 5.099503077952702
$bf[4,2] = 10.81145619092493
$bf[4,3] = 11.12552560927621
$bf[4,4] = 53.78832255196443
$bf[5,0] = 17.32067381202976
$bf[5,1] = 5.125745775879441
$bf[5,2] = 10.85382542177832
$bf[5,3] = 11.12941144369368
$bf[5,4] = 54.12458345809917
$bf[6,0] = 17.37429161009148
$bf[6,1] = 5.236457678353255
$bf[6,2] = 11.03648823692757
$bf[6,3] = 11.14969622488087
$bf[6,4] = 55.58782715190917
$bf[7,0] = 17.53178714232591
$bf[7,1] = 5.435214008443506
$bf[7,2] = 11.37827846676515
$bf[7,3] = 11.20065213200512
$bf[7,4] = 58.36509350137226
$bf[8,0] = 17.67869170758843
$bf[8,1] = 5.570033350691983
$bf[8,2] = 11.61858511635895
$bf[8,3] = 11.24480284124079
$bf[8,4] = 60.33429461897334
$bf[9,0] = 17.7520943727763
$bf[9,1] = 5.628931918819961
$bf[9,2] = 11.7253984810512
$bf[9,3] = 11.26632042125931
$bf[9,4] = 61.21173445020478
$bf[10,0] = 17.78081536612216
$bf[10,1] = 5.650885577633562
$bf[10,2] = 11.76547277334707
$bf[10,3] = 11.2746721117265
$bf[10,4] = 61.54113215065974
$bf[11,0] = 17.77458899381421
$bf[11,1] = 5.646173031961111
$bf[11,2] = 11.75685888559399
$bf[11,3] = 11.27286442143383
$bf[11,4] = 61.47032159145448
$bf[12,0] = 17.75443885218401
$bf[12,1] = 5.630745074331934
$bf[12,2] = 11.72870298118655
$bf[12,3] = 11.2670034582536
$bf[12,4] = 61.238892973049
$bf[13,0] = 17.74221612893112
$bf[13,1] = 5.621249424181437
$bf[13,2] = 11.7114076256695
$bf[13,3] = 11.2634398624792
$bf[13,4] = 61.09675597878292
$bf[14,0] = 17.67402527167919
$bf[14,1] = 5.566135233863962
$bf[14,2] = 11.61155306125455
$bf[14,3] = 11.24342523505691
$bf[14,4] = 60.27656245110112
$bf[15,0] = 17.63386247659766
$bf[15,1] = 5.531701444363895
$bf[15,2] = 11.54964449624119
$bf[15,3] = 11.2315121691926
$bf[15,4] = 59.76852985596709
$bf[16,0] = 17.61138213081525
$bf[16,1] = 5.511667190922333
$bf[16,2] = 11.51380148566128
$bf[16,3] = 11.2247951593154
$bf[16,4] = 59.47460764524109
$bf[17,0] = 17.60387778875306
$bf[17,1] = 5.504844644128221
$bf[17,2] = 11.50162574763246
$bf[17,3] = 11.22254417095536
$bf[17,4] = 59.3748032707559
$bf[18,0] = 17.63807382518238
$bf[18,1] = 5.535390667737983
$bf[18,2] = 11.55625917422693
$bf[18,3] = 11.23276637462379
$bf[18,4] = 59.82279003051616
$bf[19,0] = 17.76033250605807
$bf[19,1] = 5.635286138781964
$bf[19,2] = 11.73698329863985
$bf[19,3] = 11.26871946567366
$bf[19,4] = 61.30694880769529
$bf[20,0] = 17.84561635706555
$bf[20,1] = 5.698533784816758
$bf[20,2] = 11.85291331348885
$bf[20,3] = 11.29340138786537
$bf[20,4] = 62.26010946330663
$bf[21,0] = 17.79961369054629
$bf[21,1] = 5.664964034080035
$bf[21,2] = 11.79124345780938
$bf[21,3] = 11.28012072316019
$bf[21,4] = 61.75300098460262
$bf[22,0] = 17.63616797362855
$bf[22,1] = 5.533723510546992
$bf[22,2] = 11.55326945892895
$bf[22,3] = 11.23219893709216
$bf[22,4] = 59.79826474258434
$bf[23,0] = 17.4836413287914
$bf[23,1] = 5.383412313841912
$bf[23,2] = 11.28765633938753
$bf[23,3] = 11.18568015730479
$bf[23,4] = 57.62524153061027
$ws.Range("B2:F25").Value = $bf

# Column H for rows 2-25
$h = New-Object "object[,]" 24,1
$h[0,0] = 7.344005520526261
$h[1,0] = 7.344005520526261
$h[2,0] = 7.344005520526261
$h[3,0] = 7.344005520526261
$h[4,0] = 7.344005520526261
$h[5,0] = 7.344005520526261
$h[6,0] = 7.344005520526261
$h[7,0] = 7.344005520526261
$h[8,0] = 7.344005520526261
$h[9,0] = 7.344005520526261
$h[10,0] = 7.344005520526261
$h[11,0] = 7.344005520526261
$h[12,0] = 7.344005520526261
$h[13,0] = 7.344005520526261
$h[14,0] = 7.344005520526261
$h[15,0] = 7.344005520526261
$h[16,0] = 7.344005520526261
$h[17,0] = 7.344005520526261
$h[18,0] = 7.344005520526261
$h[19,0] = 7.344005520526261
$h[20,0] = 7.344005520526261
$h[21,0] = 7.344005520526261
$h[22,0] = 7.344005520526261
$h[23,0] = 7.344005520526261
$ws.Range("H2:H25").Value = $h

# Columns K:L for rows 2-25
$kl = New-Object "object[,]" 24,2
$kl[0,0] = 13.68452227586695
$kl[0,1] = 10.21633361865065
$kl[1,0] = 13.68601242488031
$kl[1,1] = 10.21941548474424
$kl[2,0] = 13.69238407058912
$kl[2,1] = 10.22337281166599
$kl[3,0] = 13.69635027585512
$kl[3,1] = 10.22550463785333
$kl[4,0] = 13.69709149257985
$kl[4,1] = 10.22588997657652
$kl[5,0] = 13.69243201855081
$kl[5,1] = 10.22339946035148
$kl[6,0] = 13.68390247110799
$kl[6,1] = 10.21696765057425
$kl[7,0] = 13.71052145337272
$kl[7,1] = 10.22073673742923
$kl[8,0] = 13.75648344542707
$kl[8,1] = 10.23347723032598
$kl[9,0] = 13.78309369605083
$kl[9,1] = 10.24143052699367
$kl[10,0] = 13.79398526552781
$kl[10,1] = 10.2447513893284
$kl[11,0] = 13.79160342619206
$kl[11,1] = 10.24402245456669
$kl[12,0] = 13.78397345530562
$kl[12,1] = 10.24169754843497
$kl[13,0] = 13.7794058204418
$kl[13,1] = 10.24031369526088
$kl[14,0] = 13.75485873689937
$kl[14,1] = 10.23300078806665
$kl[15,0] = 13.74125711328747
$kl[15,1] = 10.22906643214131
$kl[16,0] = 13.73397090240668
$kl[16,1] = 10.2270066597308
$kl[17,0] = 13.73159629339996
$kl[17,1] = 10.22634418068126
$kl[18,0] = 13.7426494757273
$kl[18,1] = 10.22946423171401
$kl[19,0] = 13.78619249534285
$kl[19,1] = 10.2423720504878
$kl[20,0] = 13.81939659740351
$kl[20,1] = 10.25260915819504
$kl[21,0] = 13.80124269540066
$kl[21,1] = 10.2469810569254
$kl[22,0] = 13.74201832642981
$kl[22,1] = 10.22928375681919
$kl[23,0] = 13.6986765167932
$kl[23,1] = 10.21796399885293
$ws.Range("K2:L25").Value = $kl

Write-Output "Updated loading_percent values for case 380kV"
